# Auto-generated edit script: update 'F' column (想去人数 / want-to-go count) values
# across all four worksheets per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3117
$ws.Range("F9").Value = 7313
$ws.Range("F10").Value = 66
$ws.Range("F13").Value = 28
$ws.Range("F14").Value = 424
$ws.Range("F17").Value = 1946
$ws.Range("F18").Value = 1778
$ws.Range("F19").Value = 1074
$ws.Range("F20").Value = 22
$ws.Range("F22").Value = 1814
$ws.Range("F23").Value = 1365
$ws.Range("F24").Value = 1228
$ws.Range("F25").Value = 639
$ws.Range("F26").Value = 50
$ws.Range("F27").Value = 1119
$ws.Range("F29").Value = 116
$ws.Range("F30").Value = 526
$ws.Range("F31").Value = 127
$ws.Range("F32").Value = 66
$ws.Range("F33").Value = 2686
$ws.Range("F34").Value = 1510
$ws.Range("F35").Value = 3011
$ws.Range("F36").Value = 2178
$ws.Range("F37").Value = 141
$ws.Range("F43").Value = 373
$ws.Range("F45").Value = 513
$ws.Range("F48").Value = 747
$ws.Range("F50").Value = 109

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F17").Value = 540
$ws.Range("F20").Value = 53
$ws.Range("F24").Value = 79
$ws.Range("F32").Value = 20

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1834
$ws.Range("F8").Value = 2887
$ws.Range("F9").Value = 1127
$ws.Range("F10").Value = 1106
$ws.Range("F12").Value = 416
$ws.Range("F13").Value = 1818
$ws.Range("F14").Value = 8140

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 3117
$ws.Range("F6").Value = 1834
$ws.Range("F7").Value = 2887
$ws.Range("F8").Value = 1127
$ws.Range("F9").Value = 1106
$ws.Range("F10").Value = 66
$ws.Range("F11").Value = 416
$ws.Range("F13").Value = 28
$ws.Range("F14").Value = 424
$ws.Range("F18").Value = 1074
$ws.Range("F19").Value = 22
$ws.Range("F21").Value = 1814
$ws.Range("F22").Value = 1365
$ws.Range("F23").Value = 1228
$ws.Range("F24").Value = 639
$ws.Range("F25").Value = 50
$ws.Range("F26").Value = 1119
$ws.Range("F28").Value = 116
$ws.Range("F30").Value = 540
$ws.Range("F31").Value = 526
$ws.Range("F32").Value = 127
$ws.Range("F33").Value = 66
$ws.Range("F34").Value = 2686
$ws.Range("F35").Value = 1510
$ws.Range("F36").Value = 3011
$ws.Range("F37").Value = 2178
$ws.Range("F38").Value = 141
$ws.Range("F44").Value = 79
$ws.Range("F45").Value = 513
$ws.Range("F49").Value = 20
